$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 640061.1
$ws.Range("I9").Value = 669132.0600000001
$ws.Range("K9").Value = 669132.0600000001
$ws.Range("M9").Value = -668963.0600000001

$ws.Range("H43").Value = 3886.7778
$ws.Range("I43").Value = 3719.8
$ws.Range("J43").Value = 4095.5
$ws.Range("K43").Value = 3719.8
$ws.Range("L43").Value = 4095.5
$ws.Range("M43").Value = -3650.8
$ws.Range("N43").Value = -4233.5

$ws.Range("H64").Value = 58827188
$ws.Range("I64").Value = 3591.3572
$ws.Range("K64").Value = 3591.3572
$ws.Range("M64").Value = -3343.3572

$ws.Range("H67").Value = 58827188
$ws.Range("I67").Value = 3591.3572
$ws.Range("K67").Value = 3591.3572
$ws.Range("M67").Value = -2733.3572

$ws.Range("H80").Value = 260.0625
$ws.Range("J80").Value = 217
$ws.Range("L80").Value = 651
$ws.Range("N80").Value = -2647

$ws.Range("H83").Value = 260.0625
$ws.Range("J83").Value = 217
$ws.Range("L83").Value = 1953
$ws.Range("N83").Value = -11937

$ws.Range("H100").Value = 1086
$ws.Range("I100").Value = 944.6
$ws.Range("K100").Value = 944.6
$ws.Range("M100").Value = -403.6

$ws.Range("H111").Value = 2148.1482
$ws.Range("J111").Value = 1095.238
$ws.Range("L111").Value = 3285.714
$ws.Range("N111").Value = -9419.714

$ws.Range("H113").Value = 19609622
$ws.Range("I113").Value = 47620510
$ws.Range("J113").Value = 1997.6
$ws.Range("K113").Value = 47620510
$ws.Range("L113").Value = 1997.6
$ws.Range("M113").Value = -47617256
$ws.Range("N113").Value = -8505.6


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3322.28
$ws.Range("I32").Value = 3296.125
$ws.Range("J32").Value = 3950
$ws.Range("K32").Value = 3296.125
$ws.Range("L32").Value = 3950
$ws.Range("M32").Value = -3009.125
$ws.Range("N32").Value = -4524

$ws.Range("H122").Value = 4003.6487
$ws.Range("I122").Value = 3897.3225
$ws.Range("J122").Value = 4553
$ws.Range("K122").Value = 11691.9675
$ws.Range("L122").Value = 13659
$ws.Range("M122").Value = -9241.967500000001
$ws.Range("N122").Value = -18559

$ws.Range("H132").Value = 18522356
$ws.Range("I132").Value = 4466
$ws.Range("K132").Value = 13398
$ws.Range("M132").Value = -10868

$ws.Range("H135").Value = 69672.25
$ws.Range("I135").Value = 39997.668
$ws.Range("J135").Value = 76520.234
$ws.Range("K135").Value = 39997.668
$ws.Range("L135").Value = 76520.234
$ws.Range("M135").Value = -34927.668
$ws.Range("N135").Value = -86660.234


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 9999
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 9999
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 9999
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -10279

$ws.Range("H99").Value = 3168.3
$ws.Range("I99").Value = 2240.6667
$ws.Range("K99").Value = 2240.6667
$ws.Range("M99").Value = -742.6667000000002

$ws.Range("H134").Value = 3993.4
$ws.Range("I134").Value = 3856.9333
$ws.Range("K134").Value = 11570.7999
$ws.Range("M134").Value = -9035.7999


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2197526.5
$ws.Range("I31").Value = 3497.6758
$ws.Range("J31").Value = 6256480
$ws.Range("K31").Value = 3497.6758
$ws.Range("L31").Value = 6256480
$ws.Range("M31").Value = -3202.6758
$ws.Range("N31").Value = -6257070

$ws.Range("H34").Value = 2197526.5
$ws.Range("I34").Value = 3497.6758
$ws.Range("J34").Value = 6256480
$ws.Range("K34").Value = 3497.6758
$ws.Range("L34").Value = 6256480
$ws.Range("M34").Value = -3295.6758
$ws.Range("N34").Value = -6256884

$ws.Range("H107").Value = 2632632.5
$ws.Range("I107").Value = 5000783
$ws.Range("J107").Value = 1354.2222
$ws.Range("K107").Value = 5000783
$ws.Range("L107").Value = 1354.2222
$ws.Range("M107").Value = -4998863
$ws.Range("N107").Value = -5194.2222

$ws.Range("H122").Value = 3126.9211
$ws.Range("I122").Value = 2024.2354
$ws.Range("K122").Value = 6072.706200000001
$ws.Range("M122").Value = -3622.706200000001


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 6166.5835
$ws.Range("I25").Value = 833
$ws.Range("J25").Value = 7944.4443
$ws.Range("K25").Value = 2499
$ws.Range("L25").Value = 23833.3329
$ws.Range("M25").Value = -2330
$ws.Range("N25").Value = -24171.3329

$ws.Range("H30").Value = 6166.5835
$ws.Range("I30").Value = 833
$ws.Range("J30").Value = 7944.4443
$ws.Range("K30").Value = 2499
$ws.Range("L30").Value = 23833.3329
$ws.Range("M30").Value = -2397
$ws.Range("N30").Value = -24037.3329

$ws.Range("H34").Value = 2349.625
$ws.Range("J34").Value = 2933.3333
$ws.Range("L34").Value = 8799.999899999999
$ws.Range("N34").Value = -8967.999899999999

$ws.Range("H36").Value = 1540
$ws.Range("I36").Value = 500
$ws.Range("J36").Value = 2233.3333
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 6699.999899999999
$ws.Range("M36").Value = -1331
$ws.Range("N36").Value = -7037.999899999999

$ws.Range("H110").Value = 8666.333000000001
$ws.Range("I110").Value = 8666.333000000001
$ws.Range("K110").Value = 25998.999
$ws.Range("M110").Value = -21908.999

$ws.Range("H129").Value = 1997.5555
$ws.Range("J129").Value = 2109.6
$ws.Range("L129").Value = 6328.799999999999
$ws.Range("N129").Value = -16328.8

$ws.Range("H133").Value = 4803.8237
$ws.Range("J133").Value = 4803.8237
$ws.Range("L133").Value = 14411.4711
$ws.Range("N133").Value = -24531.4711


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3503.484
$ws.Range("I122").Value = 3088.8
$ws.Range("J122").Value = 5231.3335
$ws.Range("K122").Value = 9266.400000000001
$ws.Range("L122").Value = 15694.0005
$ws.Range("M122").Value = -6816.400000000001
$ws.Range("N122").Value = -20594.0005

$ws.Range("H132").Value = 4587.8945
$ws.Range("I132").Value = 4204.375
$ws.Range("J132").Value = 6633.3335
$ws.Range("K132").Value = 12613.125
$ws.Range("L132").Value = 19900.0005
$ws.Range("M132").Value = -10083.125
$ws.Range("N132").Value = -24960.0005


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 31250.5
$ws.Range("I61").Value = 4998.5
$ws.Range("K61").Value = 4998.5
$ws.Range("M61").Value = -4796.5

$ws.Range("H68").Value = 2357.7693
$ws.Range("I68").Value = 2697.3333
$ws.Range("J68").Value = 2066.7144
$ws.Range("K68").Value = 2697.3333
$ws.Range("L68").Value = 2066.7144
$ws.Range("M68").Value = -1948.3333
$ws.Range("N68").Value = -3564.7144

$ws.Range("H71").Value = 2357.7693
$ws.Range("I71").Value = 2697.3333
$ws.Range("J71").Value = 2066.7144
$ws.Range("K71").Value = 13486.6665
$ws.Range("L71").Value = 10333.572
$ws.Range("M71").Value = -9742.666499999999
$ws.Range("N71").Value = -17821.572

$ws.Range("H100").Value = 996.1429000000001
$ws.Range("I100").Value = 962.1667
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 962.1667
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = -421.1667
$ws.Range("N100").Value = -2282

$ws.Range("H113").Value = 31250.5
$ws.Range("I113").Value = 4998.5
$ws.Range("K113").Value = 4998.5
$ws.Range("M113").Value = -2828.5

$ws.Range("H132").Value = 2996.6667
$ws.Range("I132").Value = 2746.8333
$ws.Range("K132").Value = 8240.499899999999
$ws.Range("M132").Value = -5710.499899999999

$ws.Range("H133").Value = 105993
$ws.Range("J133").Value = 105993
$ws.Range("L133").Value = 105993
$ws.Range("N133").Value = -111053


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H103").Value = 87400.664
$ws.Range("J103").Value = 87400.664
$ws.Range("L103").Value = 87400.664
$ws.Range("N103").Value = -89744.664

$ws.Range("H107").Value = 880
$ws.Range("I107").Value = 703
$ws.Range("J107").Value = 1256.125
$ws.Range("K107").Value = 2109
$ws.Range("L107").Value = 3768.375
$ws.Range("M107").Value = -189
$ws.Range("N107").Value = -7608.375

